# Remove the duplicated NOMBRE_ARTICULO_2 column (EA), which duplicated
# the existing column DZ. Deleting the entire column shifts every
# subsequent column (EB..EO) one position to the left (EA..EN).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("EA:EA").Delete()
